$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2204.611
$ws.Range("I40").Value = 2268
$ws.Range("J40").Value = 2125.375
$ws.Range("K40").Value = 2268
$ws.Range("L40").Value = 2125.375
$ws.Range("M40").Value = -2093
$ws.Range("N40").Value = -2475.375

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H64").Value = 4755.5
$ws.Range("I64").Value = 3939.0908
$ws.Range("J64").Value = 5446.3076
$ws.Range("K64").Value = 3939.0908
$ws.Range("L64").Value = 5446.3076
$ws.Range("M64").Value = -3691.0908
$ws.Range("N64").Value = -5942.3076

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H67").Value = 4755.5
$ws.Range("I67").Value = 3939.0908
$ws.Range("J67").Value = 5446.3076
$ws.Range("K67").Value = 3939.0908
$ws.Range("L67").Value = 5446.3076
$ws.Range("M67").Value = -3081.0908
$ws.Range("N67").Value = -7162.3076

$ws.Range("H74").Value = 3386.6667
$ws.Range("I74").Value = 3079.2
$ws.Range("J74").Value = 3771
$ws.Range("K74").Value = 3079.2
$ws.Range("L74").Value = 3771
$ws.Range("M74").Value = -2143.2
$ws.Range("N74").Value = -5643

$ws.Range("H76").Value = 3177879.8
$ws.Range("I76").Value = 4276619.5
$ws.Range("K76").Value = 4276619.5
$ws.Range("M76").Value = -4276304.5

$ws.Range("H77").Value = 3386.6667
$ws.Range("I77").Value = 3079.2
$ws.Range("J77").Value = 3771
$ws.Range("K77").Value = 15396
$ws.Range("L77").Value = 18855
$ws.Range("M77").Value = -10716
$ws.Range("N77").Value = -28215

$ws.Range("H79").Value = 3177879.8
$ws.Range("I79").Value = 4276619.5
$ws.Range("K79").Value = 4276619.5
$ws.Range("M79").Value = -4275527.5

$ws.Range("H113").Value = 5863.8887
$ws.Range("I113").Value = 5185
$ws.Range("J113").Value = 6203.3335
$ws.Range("K113").Value = 5185
$ws.Range("L113").Value = 6203.3335
$ws.Range("M113").Value = -1931
$ws.Range("N113").Value = -12711.3335

$ws.Range("H133").Value = 45851.668
$ws.Range("J133").Value = 45851.668
$ws.Range("L133").Value = 45851.668
$ws.Range("N133").Value = -55971.668

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15882.122
$ws.Range("I32").Value = 3243.4458
$ws.Range("J32").Value = 85816.13
$ws.Range("K32").Value = 3243.4458
$ws.Range("L32").Value = 85816.13
$ws.Range("M32").Value = -2956.4458
$ws.Range("N32").Value = -86390.13

$ws.Range("H63").Value = 6605.905
$ws.Range("I63").Value = 6836.1113
$ws.Range("J63").Value = 6433.25
$ws.Range("K63").Value = 6836.1113
$ws.Range("L63").Value = 6433.25
$ws.Range("M63").Value = -6150.1113
$ws.Range("N63").Value = -7805.25

$ws.Range("H66").Value = 6605.905
$ws.Range("I66").Value = 6836.1113
$ws.Range("J66").Value = 6433.25
$ws.Range("K66").Value = 34180.5565
$ws.Range("L66").Value = 32166.25
$ws.Range("M66").Value = -30748.5565
$ws.Range("N66").Value = -39030.25

$ws.Range("H133").Value = 47465.25
$ws.Range("J133").Value = 47465.25
$ws.Range("L133").Value = 47465.25
$ws.Range("N133").Value = -52525.25

$ws.Range("H139").Value = 50857.5
$ws.Range("J139").Value = 50857.5
$ws.Range("L139").Value = 50857.5
$ws.Range("N139").Value = -61137.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 54450
$ws.Range("J59").Value = 54450
$ws.Range("L59").Value = 54450
$ws.Range("N59").Value = -56144

$ws.Range("H86").Value = 1430.6
$ws.Range("I86").Value = 1266.4445
$ws.Range("K86").Value = 1266.4445
$ws.Range("M86").Value = -143.4445000000001

$ws.Range("H89").Value = 1430.6
$ws.Range("I89").Value = 1266.4445
$ws.Range("K89").Value = 6332.2225
$ws.Range("M89").Value = -716.2224999999999

$ws.Range("H105").Value = 241157.64
$ws.Range("I105").Value = 2817.8215
$ws.Range("K105").Value = 2817.8215
$ws.Range("M105").Value = -1070.8215

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 17039
$ws.Range("I62").Value = 20669.084
$ws.Range("K62").Value = 20669.084
$ws.Range("M62").Value = -20045.084

$ws.Range("H65").Value = 17039
$ws.Range("I65").Value = 20669.084
$ws.Range("K65").Value = 103345.42
$ws.Range("M65").Value = -100225.42

$ws.Range("H134").Value = 2558.353
$ws.Range("I134").Value = 1035.08
$ws.Range("K134").Value = 3105.24
$ws.Range("M134").Value = -570.2399999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2993.8647
$ws.Range("I136").Value = 1433.3334
$ws.Range("J136").Value = 3131.5588
$ws.Range("K136").Value = 4300.0002
$ws.Range("L136").Value = 9394.6764
$ws.Range("M136").Value = 799.9997999999996
$ws.Range("N136").Value = -19594.6764

$ws.Range("H140").Value = 8095.485
$ws.Range("I140").Value = 13042.353
$ws.Range("K140").Value = 39127.05899999999
$ws.Range("M140").Value = -33947.05899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5253.0356
$ws.Range("I70").Value = 5386.522
$ws.Range("J70").Value = 4639
$ws.Range("K70").Value = 5386.522
$ws.Range("L70").Value = 4639
$ws.Range("M70").Value = -5116.522
$ws.Range("N70").Value = -5179

$ws.Range("H73").Value = 5253.0356
$ws.Range("I73").Value = 5386.522
$ws.Range("J73").Value = 4639
$ws.Range("K73").Value = 5386.522
$ws.Range("L73").Value = 4639
$ws.Range("M73").Value = -4450.522
$ws.Range("N73").Value = -6511

$ws.Range("H80").Value = 2962.682
$ws.Range("I80").Value = 2793.8235
$ws.Range("J80").Value = 3536.8
$ws.Range("K80").Value = 2793.8235
$ws.Range("L80").Value = 3536.8
$ws.Range("M80").Value = -1795.8235
$ws.Range("N80").Value = -5532.8

$ws.Range("H83").Value = 2962.682
$ws.Range("I83").Value = 2793.8235
$ws.Range("J83").Value = 3536.8
$ws.Range("K83").Value = 13969.1175
$ws.Range("L83").Value = 17684
$ws.Range("M83").Value = -8977.1175
$ws.Range("N83").Value = -27668

$ws.Range("H138").Value = 77204.836
$ws.Range("J138").Value = 77204.836
$ws.Range("L138").Value = 77204.836
$ws.Range("N138").Value = -87484.836

$ws.Range("H139").Value = 43061.4
$ws.Range("J139").Value = 43061.4
$ws.Range("L139").Value = 43061.4
$ws.Range("N139").Value = -53341.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H96").Value = 920
$ws.Range("I96").Value = 600
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 600
$ws.Range("L96").Value = 1000
$ws.Range("M96").Value = 773
$ws.Range("N96").Value = -3746
